# Insert a new weekly price record for "Macroferia Regional de Talca - Brócoli".
# The new record becomes row 334; every existing data row from the old 334
# down to the old 366 shifts down by one (to 335..367), so the sheet grows
# from 366 to 367 rows total (A1:R367).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 334:366 down one row, leaving a blank row 334 to fill in.
$ws.Range("A334").EntireRow.Insert()

# Populate the new row 334 with the new weekly record.
$ws.Cells.Item(334, 1).Value = 5
$ws.Cells.Item(334, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(334, 3).Value = "Maule"
$ws.Cells.Item(334, 4).Value = 44769
$ws.Cells.Item(334, 5).Value = 7
$ws.Cells.Item(334, 6).Value = 100112023
$ws.Cells.Item(334, 7).Value = "Brócoli"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 5000
$ws.Cells.Item(334, 11).Value = 800
$ws.Cells.Item(334, 12).Value = 800
$ws.Cells.Item(334, 13).Value = 800
$ws.Cells.Item(334, 14).Value = "$/unidad"
$ws.Cells.Item(334, 15).Value = "Región del Maule"
$ws.Cells.Item(334, 16).Value = 800
$ws.Cells.Item(334, 17).Value = 1
$ws.Cells.Item(334, 18).Value = "Hortaliza"
